$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet holds a simple product table (A=article code, B=description,
# C="RPOTN", D="1", E=sequence number, F=unit) ending at row 43. Two new
# price-tag rows are appended at rows 44-45, continuing the same layout
# and picking up the same bordered "General" cell style (s="1") used by
# every other data row.
#
# Row 42 is used as the formatting donor because, unlike row 43, it has
# all six columns (A:F) populated/styled - copying its format ensures the
# new F44/F45 cells also get the bordered style.

function Add-ProductRow {
    param($RowNum, $Code, $Description, $Seq)

    $srcRow = "A42:F42"
    $dstRow = "A" + $RowNum + ":F" + $RowNum

    # 1) Clone the bordered "General" format from row 42 onto the new row.
    $ws.Range($srcRow).Copy()
    $ws.Range($dstRow).PasteSpecial(-4122)

    # 2) Enter the values left-to-right. The article code, the constant
    #    "1", and the sequence number all look numeric, so each of those
    #    cells is briefly switched to Text format before the value is
    #    typed in - otherwise Excel would auto-convert them to numbers.
    $ws.Range("A" + $RowNum).NumberFormat = "@"
    $ws.Range("A" + $RowNum).Value = $Code

    $ws.Range("B" + $RowNum).Value = $Description

    $ws.Range("C" + $RowNum).Value = "RPOTN"

    $ws.Range("D" + $RowNum).NumberFormat = "@"
    $ws.Range("D" + $RowNum).Value = "1"

    $ws.Range("E" + $RowNum).NumberFormat = "@"
    $ws.Range("E" + $RowNum).Value = $Seq

    $ws.Range("F" + $RowNum).Value = "PT"

    # 3) Re-paste the bordered "General" format on top. The cells already
    #    hold text values at this point, so re-applying the format only
    #    restores the plain bordered style (matching the rest of the
    #    table) without turning the text back into numbers.
    $ws.Range($srcRow).Copy()
    $ws.Range($dstRow).PasteSpecial(-4122)
}

Add-ProductRow 44 "20118161" "PRICE TAG PUTIH N222" "14"
Add-ProductRow 45 "20118163" "PRICE TAG PROMO N338" "15"
